# Applies the scheduled-runner profit-recalculation update to the Sheets workbook.
# Updates currentAveragePrice / LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ
# columns (H, I, J, K, L, M, N) for the affected leve rows across several sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 951.5
$ws.Range("I58").Value = 457.1111
$ws.Range("J58").Value = 1587.1428
$ws.Range("K58").Value = 1371.3333
$ws.Range("L58").Value = 4761.428400000001
$ws.Range("M58").Value = -1221.3333
$ws.Range("N58").Value = -5061.428400000001

$ws.Range("H70").Value = 5667
$ws.Range("I70").Value = 2987.5
$ws.Range("K70").Value = 8962.5
$ws.Range("M70").Value = -8692.5

$ws.Range("H73").Value = 5667
$ws.Range("I73").Value = 2987.5
$ws.Range("K73").Value = 8962.5
$ws.Range("M73").Value = -8026.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 14083.333
$ws.Range("J44").Value = 14083.333
$ws.Range("L44").Value = 14083.333
$ws.Range("N44").Value = -15059.333

$ws.Range("H75").Value = 55173
$ws.Range("J75").Value = 55173
$ws.Range("L75").Value = 55173
$ws.Range("N75").Value = -56921

$ws.Range("H78").Value = 55173
$ws.Range("J78").Value = 55173
$ws.Range("L78").Value = 165519
$ws.Range("N78").Value = -174255

$ws.Range("H102").Value = 7666.1665
$ws.Range("I102").Value = 4493
$ws.Range("J102").Value = 9252.75
$ws.Range("K102").Value = 4493
$ws.Range("L102").Value = 9252.75
$ws.Range("M102").Value = -2871
$ws.Range("N102").Value = -12496.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 191.15384
$ws.Range("I22").Value = 217.63637
$ws.Range("J22").Value = 45.5
$ws.Range("K22").Value = 217.63637
$ws.Range("L22").Value = 45.5
$ws.Range("M22").Value = -44.63637
$ws.Range("N22").Value = -391.5

$ws.Range("H35").Value = 96969
$ws.Range("J35").Value = 96969
$ws.Range("L35").Value = 96969
$ws.Range("N35").Value = -97589

$ws.Range("H88").Value = 16799.334
$ws.Range("J88").Value = 16799.334
$ws.Range("L88").Value = 16799.334
$ws.Range("N88").Value = -17611.334

$ws.Range("H91").Value = 16799.334
$ws.Range("J91").Value = 16799.334
$ws.Range("L91").Value = 16799.334
$ws.Range("N91").Value = -19607.334

$ws.Range("H95").Value = 7326.6665
$ws.Range("J95").Value = 7326.6665
$ws.Range("L95").Value = 7326.6665
$ws.Range("N95").Value = -12818.6665

$ws.Range("H100").Value = 18820.5
$ws.Range("J100").Value = 18820.5
$ws.Range("L100").Value = 18820.5
$ws.Range("N100").Value = -20984.5

$ws.Range("H103").Value = 26266.666
$ws.Range("J103").Value = 24400
$ws.Range("L103").Value = 24400
$ws.Range("N103").Value = -26744

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 19061.666
$ws.Range("J28").Value = 19061.666
$ws.Range("L28").Value = 19061.666
$ws.Range("N28").Value = -19551.666

$ws.Range("H31").Value = 6227.9414
$ws.Range("I31").Value = 1984.375
$ws.Range("K31").Value = 1984.375
$ws.Range("M31").Value = -1689.375

$ws.Range("H34").Value = 6227.9414
$ws.Range("I34").Value = 1984.375
$ws.Range("K34").Value = 1984.375
$ws.Range("M34").Value = -1782.375

$ws.Range("H43").Value = 39578.5
$ws.Range("J43").Value = 39578.5
$ws.Range("L43").Value = 39578.5
$ws.Range("N43").Value = -39946.5

$ws.Range("H68").Value = 83753
$ws.Range("J68").Value = 83753
$ws.Range("L68").Value = 83753
$ws.Range("N68").Value = -85251

$ws.Range("H71").Value = 83753
$ws.Range("J71").Value = 83753
$ws.Range("L71").Value = 251259
$ws.Range("N71").Value = -258747

$ws.Range("H95").Value = 21156
$ws.Range("J95").Value = 21156
$ws.Range("L95").Value = 21156
$ws.Range("N95").Value = -26648

$ws.Range("H101").Value = 39578.5
$ws.Range("J101").Value = 39578.5
$ws.Range("L101").Value = 39578.5
$ws.Range("N101").Value = -46068.5

$ws.Range("H104").Value = 40285
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 40285
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 40285
$ws.Range("M104").Value = $null
$ws.Range("N104").Value = -45527

$ws.Range("H106").Value = 27835.5
$ws.Range("J106").Value = 27835.5
$ws.Range("L106").Value = 27835.5
$ws.Range("N106").Value = -30359.5

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").Value = $null

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").Value = $null

$ws.Range("H43").Value = 15293.077
$ws.Range("I43").Value = 940
$ws.Range("J43").Value = 19599
$ws.Range("K43").Value = 940
$ws.Range("L43").Value = 19599
$ws.Range("M43").Value = -789
$ws.Range("N43").Value = -19901

$ws.Range("H62").Value = 50000
$ws.Range("J62").Value = 50000
$ws.Range("L62").Value = 50000
$ws.Range("N62").Value = -51372

$ws.Range("H65").Value = 50000
$ws.Range("J65").Value = 50000
$ws.Range("L65").Value = 150000
$ws.Range("N65").Value = -156864

$ws.Range("H68").Value = 45267.5
$ws.Range("J68").Value = 45267
$ws.Range("L68").Value = 45267
$ws.Range("N68").Value = -46889

$ws.Range("H71").Value = 45267.5
$ws.Range("J71").Value = 45267
$ws.Range("L71").Value = 135801
$ws.Range("N71").Value = -143913

$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").Value = $null

$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").Value = $null

$ws.Range("H105").Value = 18333.334
$ws.Range("J105").Value = 18333.334
$ws.Range("L105").Value = 18333.334
$ws.Range("N105").Value = -25321.334

$ws.Range("H106").Value = 21280
$ws.Range("J106").Value = 21280
$ws.Range("L106").Value = 21280
$ws.Range("N106").Value = -23804

$ws.Range("H113").Value = 2604.4614
$ws.Range("I113").Value = 1205.8
$ws.Range("J113").Value = 7266.6665
$ws.Range("K113").Value = 1205.8
$ws.Range("L113").Value = 7266.6665
$ws.Range("M113").Value = 964.2
$ws.Range("N113").Value = -11606.6665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1995.25
$ws.Range("I22").Value = 2327
$ws.Range("K22").Value = 2327
$ws.Range("M22").Value = -2032

$ws.Range("H27").Value = 1995.25
$ws.Range("I27").Value = 2327
$ws.Range("K27").Value = 2327
$ws.Range("M27").Value = -2220

$ws.Range("H100").Value = 9374.875
$ws.Range("I100").Value = 5000
$ws.Range("K100").Value = 5000
$ws.Range("M100").Value = -4459

$ws.Range("H105").Value = 18000
$ws.Range("J105").Value = 18000
$ws.Range("L105").Value = 18000
$ws.Range("N105").Value = -24988

$ws.Range("H136").Value = 2250
$ws.Range("J136").Value = 3000
$ws.Range("L136").Value = 9000
$ws.Range("N136").Value = -14100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").Value = $null

$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").Value = $null

$ws.Range("H100").Value = 3258
$ws.Range("I100").Value = 850
$ws.Range("K100").Value = 1700
$ws.Range("M100").Value = -1159

$ws.Range("H105").Value = 28499.5
$ws.Range("J105").Value = 28499.5
$ws.Range("L105").Value = 28499.5
$ws.Range("N105").Value = -35487.5

$ws.Range("H136").Value = 2633.6
$ws.Range("I136").Value = 1880.0625
$ws.Range("J136").Value = 5647.75
$ws.Range("K136").Value = 5640.1875
$ws.Range("L136").Value = 16943.25
$ws.Range("M136").Value = -3090.1875
$ws.Range("N136").Value = -22043.25
